$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values per the repull/recalculation of data
$ws.Range("F2").Value = 4
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 4
$ws.Range("F8").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -1
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = -2
$ws.Range("F22").Value = 3
$ws.Range("F23").Value = -3
$ws.Range("F24").Value = -2
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 9
$ws.Range("F27").Value = 2
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = -3
$ws.Range("F30").Value = 2
$ws.Range("F32").Value = 2
$ws.Range("F33").Value = 0
